$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
$ws.Range("G1").Value = "Message"
$ws.Range("H1").Value = "Event"

# Bring the new H1 header cell in line with the rest of the header row
# (bold font + border) before the bulk alignment pass below so the whole
# row transitions through formatting together.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# --- New "Event" column (link-click tracking) -----------------------------
$ws.Range("H2").Value = "Link click tracked for Mary Davis Cable: opened"
$ws.Range("H3").Value = "Link click tracked for Milestone Chiropractic: opened"
$ws.Range("H4").Value = "Link click tracked for Morehead Demolition: opened"
$ws.Range("H5").Value = "Link click tracked for Motor King LLC: opened"

# --- Column widths ---------------------------------------------------------
# ColumnWidth (character units) renders to a stored width that is
# ColumnWidth + 5/6, so subtract that offset to land on the exact target.
$offset = 0.8333333333333334
$ws.Columns.Item(1).ColumnWidth = 24 - $offset
$ws.Columns.Item(2).ColumnWidth = 16 - $offset
$ws.Columns.Item(3).ColumnWidth = 16 - $offset
$ws.Columns.Item(4).ColumnWidth = 7 - $offset
$ws.Columns.Item(5).ColumnWidth = 30 - $offset
$ws.Columns.Item(6).ColumnWidth = 21 - $offset
$ws.Columns.Item(7).ColumnWidth = 30 - $offset
$ws.Columns.Item(8).ColumnWidth = 30 - $offset

# --- Wrap text formatting ---------------------------------------------------
# Header row (A1:H1): keep bold + border, drop the old center/top alignment
# in favor of wrap text only.
$headerRange = $ws.Range("A1:H1")
$headerRange.HorizontalAlignment = 1
$headerRange.VerticalAlignment = -4107
$headerRange.WrapText = $true

# Data rows (A2:H5): plain cells, wrap text only.
$dataRange = $ws.Range("A2:H5")
$dataRange.WrapText = $true
